$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 already has the exact cell-style pattern that the new row 8 needs
# (borders/number formats for name/date/amount/reason/reimbursed/reimburser/amt-reimbursed),
# so copy its formatting down to row 8 before filling in the values.
$ws.Range("B5:H5").Copy()
$ws.Range("B8:H8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new expense row (row 8)
$ws.Range("B8").Value = "Michael"
$ws.Range("C8").Value = "2/18/2014"
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = "60 pages printed"
$ws.Range("F8").Value = "No"
$ws.Range("G8").Value = "N/A"
$ws.Range("H8").Value = 0

# Update header/footer date
$ws.PageSetup.LeftHeader = "Updated February 19, 2014" + [char]10
$ws.PageSetup.RightHeader = "&P"

# Update view: scroll to row 16, select E10
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("E10").Select()
